$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 11: E11 Compression -> Pin, add F11 = power in [Pin]
$ws.Range("E11").Value = "Pin"
$ws.Range("F11").Value = "power in [Pin]"

# Update existing row 12: E12 Max Gain -> Compression
$ws.Range("E12").Value = "Compression"

# Update existing row 13: remove D13 (Pin), E13 Freq -> Max Gain, remove F13
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = "Max Gain"
$ws.Range("F13").ClearContents()

# New row 14: test 12, Pin / Freq / freq in [Freq]
$ws.Range("B14").Value = 12
$ws.Range("B14").Font.Bold = $true
$ws.Range("D14").Value = "Pin"
$ws.Range("E14").Value = "Freq"
$ws.Range("F14").Value = "freq in [Freq]"

# Update selection to F13
$ws.Range("F13").Select()

# Update window position (xWindow 930 -> 1860). The headless COM-interop
# surface has no Window object model (no Left/xWindow persistence hook),
# so this is attempted best-effort and is a no-op in this runtime.
try { $excel.ActiveWindow.Left = 1860 } catch { }
